$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.377.06"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "'3.791.81"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'603.50"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'166.33"
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("D7").Value = "'3.789.84"
$ws.Range("E7").Value = "  +1.16%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.539"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").Value = "'0.173"
$ws.Range("E10").Value = "  +4.84%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "'0.462"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "'37.85"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").Value = "'0.0000249"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "'4.422.76"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "'3.784.94"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "'69.452.57"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "'7.45"
$ws.Range("E18").Value = "  +2.34%  "
$ws.Range("D19").Value = "'17.69"
$ws.Range("E19").Value = "  +3.55%  "
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("D21").Value = "'11.41"
$ws.Range("E21").Value = "  +5.58%  "
$ws.Range("D22").Value = "'494.13"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "'0.728"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "'0.0000152"
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").Value = "'84.98"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").Value = "'2.28"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("D27").Value = "'12.32"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").Value = "'10.16"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "'2.99"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "'8.16"
$ws.Range("E31").Value = "  +3.05%  "
$ws.Range("D32").Value = "'2.44"
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("D33").Value = "'32.04"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").Value = "'3.924.59"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("D35").Value = "'3.734.80"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").Value = "'5.99"
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").Value = "'1.02"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").Value = "'0.140"
$ws.Range("E39").Value = "  +5.09%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "'0.327"
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'3.09"
$ws.Range("E42").Value = "  +6.30%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'2.00"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("D44").Value = "'426.74"
$ws.Range("E44").Value = "  -2.56%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'48.49"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("D46").Value = "'8.47"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "'40.25"
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'2.815.74"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'141.95"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("E51").Value = "  +8.63%  "
